$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on A52 (the "Source:" URL row) - the new
# layout no longer marks it up as a clickable hyperlink.
$ws.Hyperlinks.Delete()

# Rewrite the "Source:" block (rows 51-54) and extend it with the "Others:"
# block (rows 55-68), finishing with the MDIPI title + long description
# (rows 71-72). Row 50 ("Source:") is untouched.
$ws.Range("A51").Value = ""
$ws.Range("A51").Style = "source"

$ws.Range("A52").Value = "Ministère du Développement industriel et de la Promotion de l'Investissement"
$ws.Range("A52").Style = "source"

$ws.Range("A53").Value = ""
$ws.Range("A53").Style = "source"

$ws.Range("A54").Value = "http://www.mdipi.gov.dz/IMG/pdf/bulletin_PME_23_francais_vf_nov_2013.pdf"
$ws.Range("A54").Style = "source"

$ws.Range("A55").Value = ""
$ws.Range("A55").Style = "source"

$ws.Range("A56").Value = "Others:"
$ws.Range("A56").Style = "source"

$ws.Range("A57").Value = ""
$ws.Range("A57").Style = "source"

$ws.Range("A58").Value = "Office National des Statistiques"
$ws.Range("A58").Style = "source"

$ws.Range("A59").Value = ""
$ws.Range("A59").Style = "source"

$ws.Range("A60").Value = "http://www.ons.dz/"
$ws.Range("A60").Style = "source"

$ws.Range("A61").Value = ""
$ws.Range("A61").Style = "source"

$ws.Range("A62").Value = "Ministere du Developpement Industriel et de la Promotion de l'Investissement"
$ws.Range("A62").Style = "source"

$ws.Range("A63").Value = ""
$ws.Range("A63").Style = "source"

$ws.Range("A64").Value = "http://www.mdipi.gov.dz/"
$ws.Range("A64").Style = "source"

$ws.Range("A65").Value = ""
$ws.Range("A65").Style = "source"

$ws.Range("A66").Value = "Agence Nationale de Developpement de la PME"
$ws.Range("A66").Style = "source"

$ws.Range("A67").Value = ""
$ws.Range("A67").Style = "source"

$ws.Range("A68").Value = "http://www.andpme.org.dz/index.php?option=com_content&view=article&id=116&Itemid=1&lang=fr"
$ws.Range("A68").Style = "source"

$ws.Range("A71").Value = "MDIPI"
$ws.Range("A71").Style = "title"

$ws.Range("A72").Value = "Ministère de l'Industrie de la Petite et Moyenne Entreprise et de la Promotion de l'investissement, Direction Générale de la Veille Stratégique, des Etudes Economiques s s et dee e et des Statistique, Bulletin de Veille, 25/02/2012, Pas de définition universelle de la PME. Available at http://www.mdipi.gov.dz/IMG/pdf/BV_20PME_20No1.pdf"
$ws.Range("A72").Style = "source"
